$d = $word.ActiveDocument

function Find-Range($searchText) {
    $rng = $d.Content.Duplicate
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $searchText"
        return $null
    }
    return $d.Range($rng.Start, $rng.End)
}

function Add-Bookmark($name, $searchText) {
    $r = Find-Range($searchText)
    $d.Bookmarks.Add($name, $r) | Out-Null
}

function Highlight-Text($searchText) {
    $rng = $d.Content.Duplicate
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Highlight = $true
    $result = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
    return $result
}

# --- Paragraph 1: "殘缺家庭或不完整的家庭" sentence ---
# Original run text:
#   一般是指核心家庭原有配偶中有一方死亡或離去,或是父母雙亡的未婚子女。(2)
# Target run split:
#   "一般是指" | [bm _Hlk65956965] "核心家庭原有配偶中有一方死亡或離去" [/bm] |
#   (highlight) ",或是" | [bm _Hlk65956986] "父母雙亡的未婚子女" [/bm] | "。(2)"

# Add the bookmarks first (while the text is still a single run) so that the
# run-merging that happens on save cannot swallow the boundaries we need.
Add-Bookmark "_Hlk65956965" "核心家庭原有配偶中有一方死亡或離去"
Add-Bookmark "_Hlk65956986" "父母雙亡的未婚子女"

# Now apply the yellow highlight to the three runs that need it.
Highlight-Text "核心家庭原有配偶中有一方死亡或離去" | Out-Null
Highlight-Text ",或是" | Out-Null
Highlight-Text "父母雙亡的未婚子女" | Out-Null

# --- Paragraph 2: "聯合家庭" sentence ---
# Original run text:
#   ,指有一個以上多核心的家庭,一般是子女成婚後繼續和父母在一個家庭裏生活,成了兩代重叠多核心家庭,或成婚後的兄弟不分家構成同胞多核心家庭,在中國通常稱作"大家庭”。
# Target run split:
#   ",指有一個以上多核心的家庭,一般是子女成婚後繼續和父母在一個家庭裏生活,成了兩代重叠多核心家庭,或" |
#   [bm _Hlk65957465] "成婚後的兄弟不分家" [/bm] |
#   "構成同胞多核心家庭,在中國通常稱作"大家庭”。"
Add-Bookmark "_Hlk65957465" "成婚後的兄弟不分家"

Write-Output "done"
